$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43; existing rows 43-57 shift down to 44-58.
$ws.Rows(43).Insert()

# Populate the newly inserted row 43 with the new data record.
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value = "Arica y Parinacota"
$ws.Range("D43").Value = 44917
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100103
$ws.Range("H43").Value = "Frutos de hueso (carozo)"
$ws.Range("I43").Value = 100103004
$ws.Range("J43").Value = "Durazno"
$ws.Range("K43").Value = "Polar King"
$ws.Range("L43").Value = "Segunda"
$ws.Range("M43").Value = 350
$ws.Range("N43").Value = 20000
$ws.Range("O43").Value = 21000
$ws.Range("P43").Value = 20429
$ws.Range("Q43").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R43").Value = "Región de Coquimbo"
$ws.Range("S43").Value = 1135
$ws.Range("T43").Value = 18
